$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5305979251861572
$ws.Range("B1").Value = 3.331262111663818
$ws.Range("C1").Value = 5.809379577636719
$ws.Range("D1").Value = 1.478864669799805
$ws.Range("E1").Value = 0.8644530773162842
